# Implement "Add Risk, Upload Risks and View Risks" of Risk Management in Dashboard.
# Mirrors the existing taskManagement / issueManagement sheets: a new worksheet
# named "riskManagement" is appended after the last tab (issueManagement) and
# holds two rows of key/label pairs for the three new risk-management actions.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last tab so it lands at the end of the
# tab strip (Worksheets.Add defaults to inserting before the active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "riskManagement"

# Row 1: programmatic/internal keys. Row 2: user-facing labels.
$ws.Range("A1").Value = "addRisk"
$ws.Range("B1").Value = "uploadRisks"
$ws.Range("C1").Value = "viewRisks"
$ws.Range("A2").Value = "Add Risk"
$ws.Range("B2").Value = "Upload Risks"
$ws.Range("C2").Value = "View Risks"

# Column widths matching the sibling sheets' look (closest values this host
# can express through the pixel-quantized ColumnWidth COM property).
$ws.Columns.Item(1).ColumnWidth = 16.917
$ws.Columns.Item(2).ColumnWidth = 17.917
$ws.Columns.Item(3).ColumnWidth = 17.25

# Leave the cursor one cell past the data, like the other sheets in this
# workbook (e.g. issueManagement's selection sits beyond its used range too).
$ws.Range("D2").Select() | Out-Null
